# Update figures in the February 2024 data workbook.
# All target cells hold numeric-looking values that are stored as *text*
# (t="inlineStr") in the original file, so we re-enter them with a leading
# apostrophe (quote-prefix) to force Excel to keep them as text rather than
# silently re-typing the cell as a number.

$wb = $excel.ActiveWorkbook

# --- Sheet "部门情况202402" (department overview) ---
$wsDept = $wb.Worksheets.Item("部门情况202402")

$wsDept.Range("K2").Value = "'82000.00"
$wsDept.Range("K3").Value = "'257000.00"
$wsDept.Range("K4").Value = "'75000.00"
$wsDept.Range("K5").Value = "'167000.00"
$wsDept.Range("K6").Value = "'84000.00"

# --- Sheet "经办人情况202402" (agent / handler overview) ---
$wsAgent = $wb.Worksheets.Item("经办人情况202402")

$wsAgent.Range("C3").Value = "'34000.00"
$wsAgent.Range("G3").Value = "'34000.00"

$wsAgent.Range("C6").Value = "'35000.00"
$wsAgent.Range("D6").Value = "'4.00"
$wsAgent.Range("G6").Value = "'35000.00"
$wsAgent.Range("H6").Value = "'4.00"

$wsAgent.Range("C8").Value = "'88000.00"
$wsAgent.Range("G8").Value = "'40000.00"

$wsAgent.Range("C10").Value = "'288004.56"
$wsAgent.Range("D10").Value = "'39.00"
$wsAgent.Range("G10").Value = "'217000.00"
$wsAgent.Range("H10").Value = "'12.00"

$wsAgent.Range("C11").Value = "'36800.00"
$wsAgent.Range("D11").Value = "'2.00"
$wsAgent.Range("G11").Value = "'18000.00"
$wsAgent.Range("H11").Value = "'1.00"

$wsAgent.Range("C14").Value = "'95569.00"
$wsAgent.Range("D14").Value = "'36.00"
$wsAgent.Range("G14").Value = "'27000.00"
$wsAgent.Range("H14").Value = "'2.00"

$wsAgent.Range("C15").Value = "'38500.00"
$wsAgent.Range("D15").Value = "'4.00"
$wsAgent.Range("G15").Value = "'29000.00"
$wsAgent.Range("H15").Value = "'3.00"

$wsAgent.Range("C17").Value = "'120500.00"
$wsAgent.Range("G17").Value = "'108000.00"

$wsAgent.Range("C19").Value = "'62000.00"
$wsAgent.Range("G19").Value = "'62000.00"
